# Update "剩余" (E column) and "开始时间" (F column) values on Sheet1.
# Rule observed from the source data refresh:
#   - For each data row (2..99) except row 36 (which has a malformed date
#     and is left untouched by the source process):
#       * if the current E (剩余) value is 1, it rolls over to 10 and the
#         F (开始时间) date advances from 20251109 to 20251119
#       * otherwise E simply decreases by 1, F stays the same

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }

    $eCell = $ws.Cells.Item($row, 5)   # column E
    $fCell = $ws.Cells.Item($row, 6)   # column F

    $eValue = $eCell.Value2

    if ($eValue -eq $null) {
        continue
    }

    if ($eValue -eq 1) {
        $eCell.Value2 = 10
        $fCell.Value2 = 20251119
    } else {
        $eCell.Value2 = $eValue - 1
    }
}
